$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Runmode column (D) fixups: tests that were "N" are now "Y"
$ws.Cells.Item(24, 4).Value = "Y"
$ws.Cells.Item(26, 4).Value = "Y"
$ws.Cells.Item(27, 4).Value = "Y"
$ws.Cells.Item(28, 4).Value = "Y"
$ws.Cells.Item(29, 4).Value = "Y"
$ws.Cells.Item(30, 4).Value = "Y"

# Remove the green highlight fill from column C (rows 2-23, 25)
$rngGreen = $ws.Range("C2:C23")
$rngGreen.Interior.ColorIndex = -4142
$rngGreen.Interior.Pattern = -4142
$ws.Cells.Item(25, 3).Interior.ColorIndex = -4142
$ws.Cells.Item(25, 3).Interior.Pattern = -4142

# Scroll position / selection reset
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A2").Select()
